$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ID[product_ID]"
$ws.Range("B1").Value = "Quantity[quantity]"
$ws.Range("C1").Value = "ProductTitle[title]"
$ws.Range("D1").Value = "UnitPrice[price]"
$ws.Range("E1").Value = "validFrom[validFrom]"
$ws.Range("F1").Value = "timestamp[timestamp]"
$ws.Range("G1").Value = "date[date]"
$ws.Range("H1").Value = "time[time]"
